$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (Exhibition) -- column F ("想去人数" / interested-count) updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 239
$ws1.Range("F3").Value = 1382
$ws1.Range("F4").Value = 21
$ws1.Range("F5").Value = 889
$ws1.Range("F7").Value = 1222
$ws1.Range("F8").Value = 1546
$ws1.Range("F10").Value = 54
$ws1.Range("F11").Value = 1975
$ws1.Range("F12").Value = 445
$ws1.Range("F13").Value = 109
$ws1.Range("F16").Value = 87
$ws1.Range("F17").Value = 80
$ws1.Range("F18").Value = 6084
$ws1.Range("F19").Value = 47
$ws1.Range("F20").Value = 5932
$ws1.Range("F21").Value = 9899
$ws1.Range("F24").Value = 181
$ws1.Range("F25").Value = 273
$ws1.Range("F26").Value = 494
$ws1.Range("F27").Value = 162
$ws1.Range("F28").Value = 145
$ws1.Range("F29").Value = 4382
$ws1.Range("F30").Value = 370

# Sheet 4: "全部类型" (All types) -- column F ("想去人数" / interested-count) updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 239
$ws4.Range("F5").Value = 1382
$ws4.Range("F6").Value = 21
$ws4.Range("F8").Value = 889
$ws4.Range("F10").Value = 1222
$ws4.Range("F12").Value = 1546
$ws4.Range("F15").Value = 1981
$ws4.Range("F17").Value = 445
$ws4.Range("F18").Value = 109
$ws4.Range("F22").Value = 87
$ws4.Range("F23").Value = 80
$ws4.Range("F24").Value = 6084
$ws4.Range("F25").Value = 47
$ws4.Range("F26").Value = 5932
$ws4.Range("F27").Value = 9900
$ws4.Range("F31").Value = 181
$ws4.Range("F32").Value = 273
$ws4.Range("F34").Value = 494
$ws4.Range("F38").Value = 162
$ws4.Range("F39").Value = 145
$ws4.Range("F40").Value = 4382
$ws4.Range("F46").Value = 370
